$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.018980872966287
$ws.Cells.Item(2, 4).Value = 1.02446114319694
$ws.Cells.Item(2, 5).Value = 0.9926147277508489
$ws.Cells.Item(2, 6).Value = 1.030335168224708
$ws.Cells.Item(2, 9).Value = 1.029215397759858
$ws.Cells.Item(2, 10).Value = 1.024186229349376
$ws.Cells.Item(2, 11).Value = 1.027289577123887
$ws.Cells.Item(2, 12).Value = 0.9955398523336033
$ws.Cells.Item(2, 13).Value = 1.033146478001981
$ws.Cells.Item(2, 14).Value = 1.025640691270855

# Row 3
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.019913273533354
$ws.Cells.Item(3, 4).Value = 1.02513220951308
$ws.Cells.Item(3, 5).Value = 0.9936372048519304
$ws.Cells.Item(3, 6).Value = 1.031500539460234
$ws.Cells.Item(3, 9).Value = 1.029395916392016
$ws.Cells.Item(3, 10).Value = 1.024755276268684
$ws.Cells.Item(3, 11).Value = 1.027768109657025
$ws.Cells.Item(3, 12).Value = 0.9963617723202692
$ws.Cells.Item(3, 13).Value = 1.034119223370288
$ws.Cells.Item(3, 14).Value = 1.026210546302058

# Row 4
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.020516462159065
$ws.Cells.Item(4, 4).Value = 1.025565968343524
$ws.Cells.Item(4, 5).Value = 0.9942998659930995
$ws.Cells.Item(4, 6).Value = 1.032254491270626
$ws.Cells.Item(4, 9).Value = 1.029510801833604
$ws.Cells.Item(4, 10).Value = 1.025122764841613
$ws.Cells.Item(4, 11).Value = 1.028076609206849
$ws.Cells.Item(4, 12).Value = 0.9968940712668345
$ws.Cells.Item(4, 13).Value = 1.034747965024949
$ws.Cells.Item(4, 14).Value = 1.026578556750913

# Row 5
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.020770009468696
$ws.Cells.Item(5, 4).Value = 1.025748207948639
$ws.Cells.Item(5, 5).Value = 0.9945786998346017
$ws.Cells.Item(5, 6).Value = 1.032571423921618
$ws.Cells.Item(5, 9).Value = 1.02955863900478
$ws.Cells.Item(5, 10).Value = 1.025277083452268
$ws.Cells.Item(5, 11).Value = 1.028206028070511
$ws.Cells.Item(5, 12).Value = 0.997117960005301
$ws.Cells.Item(5, 13).Value = 1.035012123047279
$ws.Cells.Item(5, 14).Value = 1.026733094511695

# Row 6
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.020812579219242
$ws.Cells.Item(6, 4).Value = 1.025778800145932
$ws.Cells.Item(6, 5).Value = 0.9946255319796338
$ws.Cells.Item(6, 6).Value = 1.032624636635983
$ws.Cells.Item(6, 9).Value = 1.029566644048021
$ws.Cells.Item(6, 10).Value = 1.025302984031158
$ws.Cells.Item(6, 11).Value = 1.028227741958509
$ws.Cells.Item(6, 12).Value = 0.9971555583673453
$ws.Cells.Item(6, 13).Value = 1.035056466671446
$ws.Cells.Item(6, 14).Value = 1.026759031872378

# Row 7
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.020519850201564
$ws.Cells.Item(7, 4).Value = 1.025568403881866
$ws.Cells.Item(7, 5).Value = 0.9943035907982488
$ws.Cells.Item(7, 6).Value = 1.032258726252532
$ws.Cells.Item(7, 9).Value = 1.029511442846971
$ws.Cells.Item(7, 10).Value = 1.025124827537201
$ws.Cells.Item(7, 11).Value = 1.028078339586947
$ws.Cells.Item(7, 12).Value = 0.9968970624462087
$ws.Cells.Item(7, 13).Value = 1.034751495365445
$ws.Cells.Item(7, 14).Value = 1.026580622375765

# Row 8
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.019296010198256
$ws.Cells.Item(8, 4).Value = 1.024688029045342
$ws.Cells.Item(8, 5).Value = 0.9929600610674301
$ws.Cells.Item(8, 6).Value = 1.030729036220263
$ws.Cells.Item(8, 9).Value = 1.029276802582981
$ws.Cells.Item(8, 10).Value = 1.02437869076697
$ws.Cells.Item(8, 11).Value = 1.027451535674335
$ws.Cells.Item(8, 12).Value = 0.995817528259106
$ws.Cells.Item(8, 13).Value = 1.033475364731474
$ws.Cells.Item(8, 14).Value = 1.025833426005738

# Row 9
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.017138417044893
$ws.Cells.Item(9, 4).Value = 1.023133173525957
$ws.Cells.Item(9, 5).Value = 0.9906006454969559
$ws.Cells.Item(9, 6).Value = 1.028032581264317
$ws.Cells.Item(9, 9).Value = 1.028848636154837
$ws.Cells.Item(9, 10).Value = 1.023058389517736
$ws.Cells.Item(9, 11).Value = 1.026338303870982
$ws.Cells.Item(9, 12).Value = 0.993918800172444
$ws.Cells.Item(9, 13).Value = 1.031221370737172
$ws.Cells.Item(9, 14).Value = 1.024511249777289

# Row 10
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.015699349548485
$ws.Cells.Item(10, 4).Value = 1.022094294984481
$ws.Cells.Item(10, 5).Value = 0.989033133672735
$ws.Cells.Item(10, 6).Value = 1.0262342816595
$ws.Cells.Item(10, 9).Value = 1.028553335697396
$ws.Cells.Item(10, 10).Value = 1.022174511273745
$ws.Cells.Item(10, 11).Value = 1.025590323998051
$ws.Cells.Item(10, 12).Value = 0.9926553831429383
$ws.Cells.Item(10, 13).Value = 1.029715138738248
$ws.Cells.Item(10, 14).Value = 1.023626116324808

# Row 11
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.015076061498879
$ws.Cells.Item(11, 4).Value = 1.021643913735907
$ws.Cells.Item(11, 5).Value = 0.988355674866747
$ws.Cells.Item(11, 6).Value = 1.025455433578786
$ws.Cells.Item(11, 9).Value = 1.028423134962715
$ws.Cells.Item(11, 10).Value = 1.021790915513347
$ws.Cells.Item(11, 11).Value = 1.025265066546617
$ws.Cells.Item(11, 12).Value = 0.9921088820399291
$ws.Cells.Item(11, 13).Value = 1.02906207288061
$ws.Cells.Item(11, 14).Value = 1.023241975814432

# Row 12
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.01484452034675
$ws.Cells.Item(12, 4).Value = 1.021476541869188
$ws.Cells.Item(12, 5).Value = 0.9881042295826724
$ws.Cells.Item(12, 6).Value = 1.025166108048748
$ws.Cells.Item(12, 9).Value = 1.028374422370938
$ws.Cells.Item(12, 10).Value = 1.021648300528998
$ws.Cells.Item(12, 11).Value = 1.025144045087072
$ws.Cells.Item(12, 12).Value = 0.9919059725120875
$ws.Cells.Item(12, 13).Value = 1.028819365823694
$ws.Cells.Item(12, 14).Value = 1.023099158300447

# Row 13
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.014894187803148
$ws.Cells.Item(13, 4).Value = 1.021512447318299
$ws.Cells.Item(13, 5).Value = 0.9881581567098651
$ws.Cells.Item(13, 6).Value = 1.025228170576452
$ws.Cells.Item(13, 9).Value = 1.028384887230554
$ws.Cells.Item(13, 10).Value = 1.021678897824112
$ws.Cells.Item(13, 11).Value = 1.02517001393039
$ws.Cells.Item(13, 12).Value = 0.9919494934313052
$ws.Cells.Item(13, 13).Value = 1.02887143317131
$ws.Cells.Item(13, 14).Value = 1.02312979904723

# Row 14
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.015056922728106
$ws.Cells.Item(14, 4).Value = 1.021630080356393
$ws.Cells.Item(14, 5).Value = 0.9883348863814464
$ws.Cells.Item(14, 6).Value = 1.025431518388797
$ws.Cells.Item(14, 9).Value = 1.028419115507397
$ws.Cells.Item(14, 10).Value = 1.021779129573084
$ws.Cells.Item(14, 11).Value = 1.02525506708153
$ws.Cells.Item(14, 12).Value = 0.9920921077337197
$ws.Cells.Item(14, 13).Value = 1.029042013280546
$ws.Cells.Item(14, 14).Value = 1.023230173136782

# Row 15
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.015157185882084
$ws.Cells.Item(15, 4).Value = 1.021702547350915
$ws.Cells.Item(15, 5).Value = 0.9884438009545853
$ws.Cells.Item(15, 6).Value = 1.025556804129788
$ws.Cells.Item(15, 9).Value = 1.028440158278932
$ws.Cells.Item(15, 10).Value = 1.02184086838804
$ws.Cells.Item(15, 11).Value = 1.025307443806764
$ws.Cells.Item(15, 12).Value = 0.9921799884222134
$ws.Cells.Item(15, 13).Value = 1.029147096164651
$ws.Cells.Item(15, 14).Value = 1.023291999627937

# Row 16
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.015740711718097
$ws.Cells.Item(16, 4).Value = 1.022124174005613
$ws.Cells.Item(16, 5).Value = 0.9890781214508737
$ws.Cells.Item(16, 6).Value = 1.026285967542707
$ws.Cells.Item(16, 9).Value = 1.028561927554817
$ws.Cells.Item(16, 10).Value = 1.022199950948831
$ws.Cells.Item(16, 11).Value = 1.025611881294649
$ws.Cells.Item(16, 12).Value = 0.9926916645766087
$ws.Cells.Item(16, 13).Value = 1.029758462460511
$ws.Cells.Item(16, 14).Value = 1.023651592127151

# Row 17
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.016106698720142
$ws.Cells.Item(17, 4).Value = 1.022388505342687
$ws.Cells.Item(17, 5).Value = 0.989476357848556
$ws.Cells.Item(17, 6).Value = 1.02674330577386
$ws.Cells.Item(17, 9).Value = 1.028637685672038
$ws.Cells.Item(17, 10).Value = 1.022424961115562
$ws.Cells.Item(17, 11).Value = 1.025802478557324
$ws.Cells.Item(17, 12).Value = 0.9930127773699352
$ws.Cells.Item(17, 13).Value = 1.030141726471555
$ws.Cells.Item(17, 14).Value = 1.023876921834128

# Row 18
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.016320157073171
$ws.Cells.Item(18, 4).Value = 1.022542633178449
$ws.Cells.Item(18, 5).Value = 0.9897087662937556
$ws.Cells.Item(18, 6).Value = 1.027010046835683
$ws.Cells.Item(18, 9).Value = 1.028681648908666
$ws.Cells.Item(18, 10).Value = 1.022556121709899
$ws.Cells.Item(18, 11).Value = 1.025913517796589
$ws.Cells.Item(18, 12).Value = 0.9932001317071769
$ws.Cells.Item(18, 13).Value = 1.030365195185839
$ws.Cells.Item(18, 14).Value = 1.024008268691554

# Row 19
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.016392938191606
$ws.Cells.Item(19, 4).Value = 1.0225951779122
$ws.Cells.Item(19, 5).Value = 0.9897880325774034
$ws.Cells.Item(19, 6).Value = 1.027100995853497
$ws.Cells.Item(19, 9).Value = 1.028696601028573
$ws.Cells.Item(19, 10).Value = 1.022600829819062
$ws.Cells.Item(19, 11).Value = 1.025951356745268
$ws.Cells.Item(19, 12).Value = 0.9932640239640975
$ws.Cells.Item(19, 13).Value = 1.030441378230237
$ws.Cells.Item(19, 14).Value = 1.024053040291359

# Row 20
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.016067433398297
$ws.Cells.Item(20, 4).Value = 1.022360150480432
$ws.Cells.Item(20, 5).Value = 0.9894336180360679
$ws.Cells.Item(20, 6).Value = 1.026694239431203
$ws.Cells.Item(20, 9).Value = 1.02862958083291
$ws.Cells.Item(20, 10).Value = 1.022400828342932
$ws.Cells.Item(20, 11).Value = 1.025782043019228
$ws.Cells.Item(20, 12).Value = 0.9929783193494215
$ws.Cells.Item(20, 13).Value = 1.030100614427675
$ws.Cells.Item(20, 14).Value = 1.023852754790193

# Row 21
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.015009002015984
$ws.Cells.Item(21, 4).Value = 1.021595442565757
$ws.Cells.Item(21, 5).Value = 0.9882828385668249
$ws.Cells.Item(21, 6).Value = 1.025371638268866
$ws.Cells.Item(21, 9).Value = 1.028409045798439
$ws.Cells.Item(21, 10).Value = 1.021749617419081
$ws.Cells.Item(21, 11).Value = 1.025230026740556
$ws.Cells.Item(21, 12).Value = 0.9920501090198102
$ws.Cells.Item(21, 13).Value = 1.028991785263843
$ws.Cells.Item(21, 14).Value = 1.023200619072136

# Row 22
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.014343384155811
$ws.Cells.Item(22, 4).Value = 1.021114176675687
$ws.Cells.Item(22, 5).Value = 0.9875604150241495
$ws.Cells.Item(22, 6).Value = 1.024539911711442
$ws.Cells.Item(22, 9).Value = 1.028268360456887
$ws.Cells.Item(22, 10).Value = 1.021339420759355
$ws.Cells.Item(22, 11).Value = 1.02488175852938
$ws.Cells.Item(22, 12).Value = 0.9914670000341481
$ws.Cells.Item(22, 13).Value = 1.028293871905317
$ws.Cells.Item(22, 14).Value = 1.022789839886102

# Row 23
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.01469625399651
$ws.Cells.Item(23, 4).Value = 1.021369348480014
$ws.Cells.Item(23, 5).Value = 0.9879432794643023
$ws.Cells.Item(23, 6).Value = 1.024980840617256
$ws.Cells.Item(23, 9).Value = 1.028343132358319
$ws.Cells.Item(23, 10).Value = 1.021556945167888
$ws.Cells.Item(23, 11).Value = 1.025066495060885
$ws.Cells.Item(23, 12).Value = 0.991776070289318
$ws.Cells.Item(23, 13).Value = 1.028663920093881
$ws.Cells.Item(23, 14).Value = 1.023007673204246

# Row 24
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.016085175755397
$ws.Cells.Item(24, 4).Value = 1.022372962983364
$ws.Cells.Item(24, 5).Value = 0.9894529299347244
$ws.Cells.Item(24, 6).Value = 1.026716410450341
$ws.Cells.Item(24, 9).Value = 1.028633243756631
$ws.Cells.Item(24, 10).Value = 1.022411733163688
$ws.Cells.Item(24, 11).Value = 1.025791277369939
$ws.Cells.Item(24, 12).Value = 0.9929938892766442
$ws.Cells.Item(24, 13).Value = 1.030119191445907
$ws.Cells.Item(24, 14).Value = 1.023863675097046

# Row 25
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.017696326556693
$ws.Cells.Item(25, 4).Value = 1.02353555214751
$ws.Cells.Item(25, 5).Value = 0.9912096547607049
$ws.Cells.Item(25, 6).Value = 1.028729794337557
$ws.Cells.Item(25, 9).Value = 1.028961066645036
$ws.Cells.Item(25, 10).Value = 1.023400369315166
$ws.Cells.Item(25, 11).Value = 1.026627130601935
$ws.Cells.Item(25, 12).Value = 0.9944092447426414
$ws.Cells.Item(25, 13).Value = 1.031804710314264
$ws.Cells.Item(25, 14).Value = 1.024853715225256

